# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E26) currently lists periods in
# descending order (2209 .. 2111). The update re-sorts that block into
# ascending order (2111 .. 2209), and the "Valor Mora" figure tied to
# period 2209 (31495, all other periods being 36341) moves along with it
# from row 16 down to row 26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("2111", "2112", "2201", "2202", "2203", "2204", "2205", "2206", "2207", "2208", "2209")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    if ($periods[$i] -eq "2209") {
        $ws.Range("F$row").Value = 31495
    } else {
        $ws.Range("F$row").Value = 36341
    }
}
